$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Rows 2-5: set column A to the practice email and column B to the password,
# clearing the fill style that those cells previously had (plain/no style).
$ws.Range("A2:B5").ClearFormats()

$ws.Range("A2:A5").Value = "tpractice427@gmail.com"
$ws.Range("B2:B5").Value = "Testing0@"

# Update the active selection on the sheet to match the saved view.
$ws.Range("B2:B5").Select()
